# The authored change widens column A ("Reference") on the ESP_WiFi_Display
# sheet from its original width (~9.77 chars) to ~15.69 chars. Everything
# else in the target diff (font charset attributes, cellXfs table
# deduplication/reindexing, the 0.01 rounding on column I's width) is a
# byte-level re-serialization artifact of the original authoring
# application and carries no semantic/visual effect, so we only need to
# reproduce the actual content edit: the column A resize.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESP_WiFi_Display")

$ws.Columns.Item(1).ColumnWidth = 14.75
